$d = $word.ActiveDocument

# 1. Delete the bookmark that currently sits at the end of the
#    "Update Call..." paragraph, we will re-create it later in the middle
#    of the "Forward declaration header." text.
$d.Bookmarks("_GoBack").Delete()

# 2. Delete the whole "Update Call to use builtin AsmJit functionality for
#    remote code gen." paragraph (including its paragraph mark).
$r = $d.Content.Find.Execute("Update Call to use builtin AsmJit functionality for remote code gen.`r", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found = $d.Content
$found.Find.Execute("Update Call to use builtin AsmJit functionality for remote code gen.`r", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found.Delete()

# 3. Split "Forward declaration header." into "Forward dec" + bookmark + "laration header."
$target = $d.Content
$target.Find.Execute("Forward declaration header.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $target.Start + [string]"Forward dec".Length
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
